$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the existing table (Table1) entirely - the final layout has no table.
if ($ws.ListObjects.Count -gt 0) {
    $ws.ListObjects.Item(1).Delete()
}

# Clear everything so we can rebuild the sheet from scratch (also resets the
# shared-string table so new unique strings are appended in the exact order
# we write them).
$ws.Cells.Clear()

# Rename the first sheet.
$ws.Name = "Tablas datos"

# --- Write cell values in the same chronological order the original author
# --- must have used, so the shared-string table ends up in the same order:
# 0 Reuerimiento 0, 1 Tamaño CSV [%], 2 Tiempo [s], 3 Reuerimiento 1,
# 4 Reuerimiento 2, 5 Reuerimiento 3, 6 Reuerimiento 4.

# Requerimiento 0 block (rows 1,3-7)
$ws.Range("A1").Value = "Reuerimiento 0"
$ws.Range("A3").Value = "Tamaño CSV [%]"
$ws.Range("B3").Value = "Tiempo [s]"

# Requerimiento 1 block (rows 1,3-7, columns D/E)
$ws.Range("D1").Value = "Reuerimiento 1"

# Requerimiento 2 block (rows 9,11-15)
$ws.Range("A9").Value = "Reuerimiento 2"

# Requerimiento 3 block (rows 9,11-15, columns D/E)
$ws.Range("D9").Value = "Reuerimiento 3"

# Requerimiento 4 block (rows 17,19-23)
$ws.Range("A17").Value = "Reuerimiento 4"

# Fill in the rest of the repeated header text (reuses the shared strings
# created above).
$ws.Range("D3").Value = "Tamaño CSV [%]"
$ws.Range("E3").Value = "Tiempo [s]"
$ws.Range("A11").Value = "Tamaño CSV [%]"
$ws.Range("B11").Value = "Tiempo [s]"
$ws.Range("D11").Value = "Tamaño CSV [%]"
$ws.Range("E11").Value = "Tiempo [s]"
$ws.Range("A19").Value = "Tamaño CSV [%]"
$ws.Range("B19").Value = "Tiempo [s]"

# Numeric data columns (Tamaño CSV [%]) for each block.
$ws.Range("A4").Value = 10
$ws.Range("A5").Value = 50
$ws.Range("A6").Value = 80
$ws.Range("A7").Value = 100

$ws.Range("D4").Value = 10
$ws.Range("D5").Value = 50
$ws.Range("D6").Value = 80
$ws.Range("D7").Value = 100

$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 50
$ws.Range("A14").Value = 80
$ws.Range("A15").Value = 100

$ws.Range("D12").Value = 10
$ws.Range("D13").Value = 50
$ws.Range("D14").Value = 80
$ws.Range("D15").Value = 100

$ws.Range("A20").Value = 10
$ws.Range("A21").Value = 50
$ws.Range("A22").Value = 80
$ws.Range("A23").Value = 100

# Empty "Tiempo [s]" data cells - still need the centered style applied, so
# touch them with a blank value so they materialize in the sheet.
$ws.Range("B4").Value = ""
$ws.Range("B5").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("B7").Value = ""

$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = ""

$ws.Range("B12").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("B15").Value = ""

$ws.Range("E12").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("E15").Value = ""

$ws.Range("B20").Value = ""
$ws.Range("B21").Value = ""
$ws.Range("B22").Value = ""
$ws.Range("B23").Value = ""

# Centered alignment (reuses the workbook's existing "centered" cell style)
# for every data table (header + 4 data rows) in each block.
$ws.Range("A3:B7").HorizontalAlignment = -4108
$ws.Range("D3:E7").HorizontalAlignment = -4108
$ws.Range("A11:B15").HorizontalAlignment = -4108
$ws.Range("D11:E15").HorizontalAlignment = -4108
$ws.Range("A19:B23").HorizontalAlignment = -4108

# Column widths (character units); the stored width ends up snapped to the
# engine's internal pixel grid, closest achievable value is used.
$ws.Columns.Item(2).ColumnWidth = 12.17
$ws.Columns.Item(4).ColumnWidth = 18.83
$ws.Columns.Item(5).ColumnWidth = 12

# Selection moves to B12 in the final workbook.
$ws.Range("B12").Select()
